# Generate Report for Handback
# Updates status/timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 16:16:56"
$wsOverview.Range("G4").Value = "2016-08-23 16:16:56"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-23 16:16:51"
$wsZhCn.Range("H4").Value = "2016-08-23 16:16:51"
$wsZhCn.Range("K3").Value = "2016-08-23 16:17:21"
$wsZhCn.Range("K4").Value = "2016-08-23 16:17:21"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 16:16:56"
$wsDeDe.Range("H4").Value = "2016-08-23 16:16:56"
$wsDeDe.Range("K3").Value = "2016-08-23 16:17:28"
$wsDeDe.Range("K4").Value = "2016-08-23 16:17:28"
